$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Customer_NIC column (B): replace bare numbers with "<NIC>V" text values
# ---------------------------------------------------------------------------
$nics = @(
    "20035301843V",
    "19977593220V",
    "20203162022V",
    "20243492022V",
    "19910751235V",
    "19900323750V",
    "19992251255V",
    "20021215135V",
    "20021151234V"
)
for ($i = 0; $i -lt $nics.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $nics[$i]
}

# ---------------------------------------------------------------------------
# 2) New Bill_amount column (D)
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Bill_amount"

$amounts = @(7800, 3225, 1350, 1698, 2185, 1235, 6598, 5396, 783)
for ($i = 0; $i -lt $amounts.Length; $i++) {
    $row = $i + 2
    $ws.Range("D$row").Value = $amounts[$i]
}

# ---------------------------------------------------------------------------
# 3) Formatting
# ---------------------------------------------------------------------------

# 3a) Column B default alignment flag (no explicit sub-values, General)
$ws.Columns.Item(2).HorizontalAlignment = 1

# 3b) B2:B10 -> right / top aligned. Built on a scratch cell + PasteSpecial so
#     the combined alignment lands on a single new style without leaving
#     unused intermediate entries behind in the style table.
$ws.Range("Z1").HorizontalAlignment = -4152
$ws.Range("Z1").VerticalAlignment = -4160
$ws.Range("Z1").Copy()
$ws.Range("B2:B10").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# 3c) Header row (A1:D1) centered
$ws.Range("A1:D1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 4) Selection cursor
# ---------------------------------------------------------------------------
$ws.Range("D19").Select()
